$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Δοκιμή")

# Row 2: refresh the sample row with a new set of source values (row 81 of the
# other sheet). B2/C2 are entered as text (apostrophe-prefixed) matching how
# they were pasted in the source edit.
$ws2.Range("A2").Value = 0.79797979797979801
$ws2.Range("B2").Value = "'1787.8315470249322"
$ws2.Range("C2").Value = "'29.103830456733704"
$ws2.Range("D2").Value = "(176.04523429126368, 187.54290076821835, 187.54290076821835, 292.0351469142523, 311.1082145879472, 311.1082145879472, 107.87620862188771, 107.9521197646562, 108.1567539274278)"
$ws2.Range("E2").Value = "(171.88042281126383, 186.2313963153644, 186.0326667717967, 289.47768714949615, 304.0576713997042, 308.3542571105595, 105.99586159281041, 106.42672757692537, 106.61534048346061)"
$ws2.Range("F2").Value = "(3382.01401871867, 3532.08445431993, 3541.24478767337, 9083.79983429210, 9737.28110427440, 9675.96869823821, 6198.74707890619, 6202.34309370489, 6219.97423860348)"

# Fix the M-column formulas to use C (not A) as the multiplier, with the
# discount rate dropping from 0.1 to 0.05 for rows 6-8.
$ws2.Range("M3").Formula = "=C3*(E3-0.1*D3)"
$ws2.Range("M4").Formula = "=C4*(E4-0.1*D4)"
$ws2.Range("M5").Formula = "=C5*(E5-0.1*D5)"
$ws2.Range("M6").Formula = "=C6*(E6-0.05*D6)"
$ws2.Range("M7").Formula = "=C7*(E7-0.05*D7)"
$ws2.Range("M8").Formula = "=C8*(E8-0.05*D8)"

$null = $ws2.Range("L7").Select()
